$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Description" column (old column B)
$ws.Columns("B").Delete()

# Apply Times New Roman font to header row and body rows
$ws.Range("A1:C12").Font.Name = "Times New Roman"

# Center-align the "Number of Features" column (B) for data rows
$ws.Range("B2:B12").HorizontalAlignment = -4108

# Bottom double-border under the last data row (row 12)
$ws.Range("A12:C12").Borders(9).LineStyle = -4119
